$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 23 with the latest processed e-mail ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A23").Value = "Sponsoraanvraag"
$logs.Range("B23").Value = "mailmind.test@zohomail.eu"
$logs.Range("C23").Value = "Zou uw bedrijf bereid zijn om ons sportevenement te sponsoren?"
$logs.Range("D23").Value = "Samenwerking / Partnerverzoek"
$logs.Range("F23").Value = "2025-06-19 21:40:10"
$logs.Range("G23").Value = "Nee"

# Extend the conditional-formatting ranges (D2:D22 -> D2:D23, G2:G22 -> G2:G23)
# so the newly added row also gets the category/answered highlighting.
$dRules = $logs.Range("D2:D22").FormatConditions
for ($i = 1; $i -le $dRules.Count; $i++) {
    $dRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D23"))
}

$gRules = $logs.Range("G2:G22").FormatConditions
for ($i = 1; $i -le $gRules.Count; $i++) {
    $gRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G23"))
}

# --- Dashboard sheet: bump the "Samenwerking / Partnerverzoek" count ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 6
